# Automatische test-sync: 2025-06-29 14:41:50
# Appends two new mail-log entries (rows 15 & 16) to the "Logs" sheet and
# updates the corresponding category counters on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Row 15: Testmail #3 (Productinformatie / verzendkosten) -----------
$logs.Range("A15").Value2 = "Wat zijn de verzendkosten?"
$logs.Range("B15").Value2 = "mailmind.test@zohomail.eu"
$logs.Range("C15").Value2 = "Testmail #3: Wat zijn de verzendkosten?"
$logs.Range("D15").Value2 = "Productinformatie"
$logs.Range("F15").Value2 = "2025-06-29 14:41:36"
$logs.Range("G15").Value2 = "Nee"
$logs.Range("H15").Value2 = "Ja"
$logs.Range("I15").Value2 = "Nee"

# --- Row 16: Testmail #1 (Openingstijden / Locatie) ---------------------
$logs.Range("A16").Value2 = "Wanneer zijn jullie open?"
$logs.Range("B16").Value2 = "mailmind.test@zohomail.eu"
$logs.Range("C16").Value2 = "Testmail #1: Wanneer zijn jullie open?"
$logs.Range("D16").Value2 = "Openingstijden / Locatie"
$logs.Range("E16").Value2 = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F16").Value2 = "2025-06-29 14:41:40"
$logs.Range("G16").Value2 = "Ja"
$logs.Range("H16").Value2 = "Nee"
$logs.Range("I16").Value2 = "Ja"

# --- Extend conditional-formatting ranges to cover the new rows ---------
foreach ($col in @("D", "G", "H", "I")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "14")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "16")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard counters --------------------------------------------------
$dashboard.Range("B2").Value2 = 5
$dashboard.Range("B4").Value2 = 4
